$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'271.94"
$ws.Range("G2").Value = "'3"
$ws.Range("G3").Value = "'3"
$ws.Range("D4").Value = "'6.256"
$ws.Range("G4").Value = "'3"
$ws.Range("D5").Value = "'0.06188"
$ws.Range("G5").Value = "'3"
$ws.Range("G6").Value = "'3"
$ws.Range("D7").Value = "'6.541"
$ws.Range("G7").Value = "'3"
$ws.Range("D8").Value = "'1.430"
$ws.Range("G8").Value = "'3"
$ws.Range("D9").Value = "'0.8233"
$ws.Range("G9").Value = "'3"
$ws.Range("D10").Value = "'0.1652"
$ws.Range("G10").Value = "'3"
$ws.Range("D11").Value = "'0.08283"
$ws.Range("G11").Value = "'3"
$ws.Range("D12").Value = "'0.03540"
$ws.Range("G12").Value = "'3"
$ws.Range("D13").Value = "'0.03181"
$ws.Range("G13").Value = "'3"
$ws.Range("D14").Value = "'0.09185"
$ws.Range("G14").Value = "'3"
$ws.Range("D15").Value = "'3.767"
$ws.Range("G15").Value = "'3"
$ws.Range("D16").Value = "'0.001632"
$ws.Range("G16").Value = "'3"
$ws.Range("D17").Value = "'0.04682"
$ws.Range("G17").Value = "'3"
$ws.Range("D18").Value = "'0.006298"
$ws.Range("G18").Value = "'3"
$ws.Range("D19").Value = "'0.006198"
$ws.Range("G19").Value = "'3"
$ws.Range("G20").Value = "'3"
$ws.Range("D21").Value = "'0.0001500"
$ws.Range("G21").Value = "'3"
$ws.Range("D22").Value = "'3.721"
$ws.Range("G22").Value = "'3"
$ws.Range("D23").Value = "'2.292"
$ws.Range("G23").Value = "'3"
$ws.Range("D24").Value = "'0.01383"
$ws.Range("G24").Value = "'3"
$ws.Range("D25").Value = "'0.3293"
$ws.Range("G25").Value = "'3"
$ws.Range("G26").Value = "'3"
$ws.Range("G27").Value = "'3"
$ws.Range("D28").Value = "'0.0002713"
$ws.Range("G28").Value = "'3"
$ws.Range("G29").Value = "'3"
$ws.Range("G30").Value = "'3"
$ws.Range("G31").Value = "'3"
$ws.Range("G32").Value = "'3"
$ws.Range("G33").Value = "'3"
$ws.Range("G34").Value = "'3"
$ws.Range("G35").Value = "'3"
$ws.Range("G36").Value = "'3"
$ws.Range("G37").Value = "'3"
$ws.Range("G38").Value = "'3"
$ws.Range("G39").Value = "'3"
$ws.Range("D40").Value = "'0.04685"
$ws.Range("G40").Value = "'3"
$ws.Range("D41").Value = "'0.007029"
$ws.Range("G41").Value = "'3"
$ws.Range("D42").Value = "'0.004600"
$ws.Range("G42").Value = "'3"
$ws.Range("D43").Value = "'0.1118"
$ws.Range("G43").Value = "'3"
$ws.Range("D44").Value = "'0.01043"
$ws.Range("G44").Value = "'3"
$ws.Range("D45").Value = "'0.00006250"
$ws.Range("G45").Value = "'3"
$ws.Range("D46").Value = "'0.0009901"
$ws.Range("G46").Value = "'3"
$ws.Range("G47").Value = "'3"
$ws.Range("D48").Value = "'0.8901"
$ws.Range("G48").Value = "'3"
$ws.Range("D49").Value = "'0.001396"
$ws.Range("G49").Value = "'3"
$ws.Range("D50").Value = "'0.00001900"
$ws.Range("G50").Value = "'3"
$ws.Range("D51").Value = "'0.01240"
$ws.Range("G51").Value = "'3"
